$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates (new addressing scheme: base network 192.168.0.0/24) ---
$ws.Range("D3").Value  = "IP 192.168.0.0/24"
$ws.Range("B6").Value  = "2^3=8 -> Tenemos 8 subredes disponibles a utilizar"
$ws.Range("B14").Value = "192.168.0.0/27"
$ws.Range("B18").Value = "2 ^ 5 - 2 -> 32"

# --- Subnet table (rows 24-28): new third-octet-per-department scheme ---
$ws.Range("D24").Value = "192.168.6.32"
$ws.Range("E24").Value = "192.168.6.33"
$ws.Range("F24").Value = "192.168.6.62"
$ws.Range("G24").Value = "192.168.6.63"

$ws.Range("D25").Value = "192.168.7.64"
$ws.Range("E25").Value = "192.168.7.65"
$ws.Range("F25").Value = "192.168.7.94"
$ws.Range("G25").Value = "192.168.7.95"
$ws.Range("H25").Value = "Departamento de Compras"

$ws.Range("D26").Value = "192.168.8.96"
$ws.Range("E26").Value = "192.168.8.97"
$ws.Range("F26").Value = "192.168.8.126"
$ws.Range("G26").Value = "192.168.8.127"
$ws.Range("H26").Value = "Departamento de Direccion General"

$ws.Range("D27").Value = "192.168.9.128"
$ws.Range("E27").Value = "192.168.9.129"
$ws.Range("F27").Value = "192.168.9.158"
$ws.Range("G27").Value = "192.168.9.159"
$ws.Range("H27").Value = "Departamento de Control de Gestión"

$ws.Range("D28").Value = "192.168.10.160"
$ws.Range("E28").Value = "192.168.10.161"
$ws.Range("F28").Value = "192.168.10.190"
$ws.Range("G28").Value = "192.168.10.191"
$ws.Range("H28").Value = "Departamento financiero"

# --- Remove the last two departments (table now only covers 6 subnets) ---
$ws.Rows("29:30").Delete()

# --- View state: scroll back to top, move selection ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("M23").Select()
